# Applies the "v1.0.2 -> v1.0.3" content changes to the "UC009 - Prestar Contas" test-suite sheet.
#
# Summary of the change:
#  1) TC3 and TC4 had their last "Steps"/"Expected Results" pair swapped:
#       TC3's last step becomes "Chefe Clica em visualizar comprovante." / "SYSTEM Exibe modal com o comprovante."
#       TC4's last step becomes "Chefe Clica para detalhar a solicitação de diária." / "SYSTEM Apresenta a tela de Detalhar Diárias"
#  2) TC7/TC8/TC9 expected-result content got rotated, and TC8 gained an extra step row
#     while TC9 lost one (it now collapses to a single step row):
#       TC7's 2nd step "Expected Results" becomes the MSG212 message.
#       TC8 gains a 2nd step whose "Expected Results" is the old TC7 message
#         ("...não está em nenhum desses dois estados...").
#       TC9's (now single) step "Expected Results" becomes the old TC8 message
#         ("...perfil OPERADOR...").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# ---------------------------------------------------------------------------
# 1) Swap TC3's and TC4's last step (Steps column B / Expected Results column D)
# ---------------------------------------------------------------------------

$tc3Step = $ws.Range("B32").Value
$tc3Expected = $ws.Range("D32").Value
$tc4Step = $ws.Range("B41").Value
$tc4Expected = $ws.Range("D41").Value

$ws.Range("B32").Value = $tc4Step
$ws.Range("D32").Value = $tc4Expected
$ws.Range("B41").Value = $tc3Step
$ws.Range("D41").Value = $tc3Expected

# ---------------------------------------------------------------------------
# 2) Rotate the TC7 / TC8 / TC9 "Expected Results" content.
# ---------------------------------------------------------------------------

# Capture the three messages (by current/original location) before touching any rows.
$msgNaoEsta = $ws.Range("D67").Value    # TC7's current 2nd-step message ("...não está em nenhum...")
$msgOperador = $ws.Range("D74").Value   # TC8's current (only) step message ("...perfil OPERADOR...")
$msgMsg212 = $ws.Range("D82").Value     # TC9's current 2nd-step message ("...MSG212...")
$genericExpected = $ws.Range("D66").Value  # generic 1st-step "Expected Results" text reused across test cases
$genericStep2 = $ws.Range("B67").Value     # generic 2nd-step "Steps" text reused across test cases

# TC7 now reports the MSG212 message on its 2nd step.
$ws.Range("D67").Value = $msgMsg212

# Remove the old TC9 2nd-step row (the one holding the MSG212 message); this shifts
# everything below it (the TC10 block) up by one row.
$ws.Rows("82:82").Delete()

# TC8's existing single step row becomes a generic first step ...
$ws.Range("D74").Value = $genericExpected

# ... and a new second step row is inserted right after it, carrying the message that
# used to belong to TC7 ("...não está em nenhum..."). Inserting here copies row 74's
# formatting into the new row 75, and shifts the TC9 block (and everything after it)
# back down by one row, restoring the original row numbering below.
$ws.Rows("75:75").Insert()
$ws.Range("A75").Value = 2
$ws.Range("B75").Value = $genericStep2
$ws.Range("C75").Value = $null
$ws.Range("D75").Value = $msgNaoEsta
$ws.Range("E75").Value = $null
$ws.Range("F75").Value = $null

# TC9 (now shifted back down to rows 77-82) collapses to a single step row: its former
# first (generic) step, at row 82, now reports the "perfil OPERADOR" message.
$ws.Range("D82").Value = $msgOperador
